$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name and title to reflect new "through" date
$ws.Name = "Through 2022-06-07"
$ws.Range("I1").Value = "2022 (through 06-07)"

# Update data values for June (row 6) and July (row 7), and Total (row 14)
$ws.Range("I6").Value = 113
$ws.Range("I7").Value = 21
$ws.Range("I14").Value = 684
